$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column L (column 12) holds values that were mistakenly scaled by 1e7.
# Divide each value in L2:L28 by 10,000,000 to correct it.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $cell.Value = $cell.Value() / 10000000
}

# Update the active selection on the sheet from N1 to N4.
$ws.Range("N4").Select()
